$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Cells whose value TYPE changes (number <-> text placeholder) ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C20").Value = 2
$ws.Range("F20").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("C26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D26").Value = 1
$ws.Range("F26").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Value = -100
$ws.Range("H26").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("C27").Value = 6
$ws.Range("F27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = "'0"
$ws.Range("A27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "'***.*"
$ws.Range("A27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("C30").Value = 1
$ws.Range("F30").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("D30").Value = 1
$ws.Range("F30").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = 0
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Plain numeric value updates (no type/style change) ---
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -7.692307692307
$ws.Range("I16").Value = 118
$ws.Range("J16").Value = 105
$ws.Range("K16").Value = 12.380952380952
$ws.Range("L16").Value = -7.8125
$ws.Range("M16").Value = 57.333333333333
$ws.Range("N16").Value = -83.791208791208
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("I17").Value = 114
$ws.Range("J17").Value = 104
$ws.Range("K17").Value = 9.615384615384
$ws.Range("L17").Value = 56.164383561643
$ws.Range("M17").Value = 93.220338983050
$ws.Range("N17").Value = -27.848101265822
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -11.538461538461
$ws.Range("I18").Value = 216
$ws.Range("J18").Value = 145
$ws.Range("K18").Value = 48.965517241379
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 44.966442953020
$ws.Range("N18").Value = -70.771312584573
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = -8
$ws.Range("F19").Value = 98
$ws.Range("G19").Value = 91
$ws.Range("H19").Value = 7.692307692307
$ws.Range("I19").Value = 1041
$ws.Range("J19").Value = 681
$ws.Range("K19").Value = 52.863436123348
$ws.Range("L19").Value = 85.231316725978
$ws.Range("M19").Value = 12.540540540540
$ws.Range("N19").Value = -69.008633521881
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 62
$ws.Range("J20").Value = 64
$ws.Range("K20").Value = -3.125
$ws.Range("L20").Value = 121.428571428571
$ws.Range("M20").Value = 77.142857142857
$ws.Range("N20").Value = -92.061459667093
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 152
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = -0.653594771241
$ws.Range("I21").Value = 1569
$ws.Range("J21").Value = 1113
$ws.Range("K21").Value = 40.970350404312
$ws.Range("L21").Value = 52.775073028237
$ws.Range("M21").Value = 25.620496397117
$ws.Range("N21").Value = -72.878133102852
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -25
$ws.Range("F22").Value = 9
$ws.Range("H22").Value = 80
$ws.Range("I22").Value = 83
$ws.Range("J22").Value = 67
$ws.Range("K22").Value = 23.880597014925
$ws.Range("L22").Value = 6.410256410256
$ws.Range("M22").Value = 50.909090909090
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = 72.727272727272
$ws.Range("F24").Value = 301
$ws.Range("G24").Value = 220
$ws.Range("H24").Value = 36.818181818181
$ws.Range("I24").Value = 3373
$ws.Range("J24").Value = 1841
$ws.Range("K24").Value = 83.215643671917
$ws.Range("L24").Value = 137.535211267606
$ws.Range("M24").Value = 132.620689655172
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = -22.5
$ws.Range("I25").Value = 318
$ws.Range("J25").Value = 306
$ws.Range("K25").Value = 3.921568627450
$ws.Range("L25").Value = 38.260869565217
$ws.Range("M25").Value = 46.543778801843
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = 66.666666666666
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 550
$ws.Range("I27").Value = 103
$ws.Range("K27").Value = 41.095890410958
$ws.Range("L27").Value = 68.852459016393
$ws.Range("I30").Value = 14
$ws.Range("J30").Value = 11
$ws.Range("K30").Value = 27.272727272727
$ws.Range("L30").Value = 250
